$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held the Melbourne exposure site; it is replaced with the
# newly reported Camberwell restaurant exposure.
$ws.Range("A2").Value = "Camberwell"
$ws.Range("B2").Value = "Tao Dumplings  1 Evans Place, Camberwell VIC 3124"
$ws.Range("C2").Value = "29/12/20 12:30pm-1:30pm"
$ws.Range("D2").Value = "Case ate at restaurant"
$ws.Range("E2").Value = "new"

# Row 3 previously held the Moorabbin exposure site; it is replaced with the
# Melbourne venue, now with a corrected exposure period.
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C3").Value = "28/12/2020 10:00pm-12.00am"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "new"

# A new row 4 is added to keep the original (now superseded) Melbourne
# exposure period on record, marked as "old".
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C4").Value = "28/12/2020 10:30pm-12.00am"
$ws.Range("D4").Value = "Case attended venue"
$ws.Range("E4").Value = "old"

$ws.Range("B4").Select()
